# RV6_swiss_travels.xlsx — rename the row label so it matches the
# dataset name used in datamanagement.py ("swiss_travels") and leave the
# selection where the author left it (A3) when they saved the file.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# A2 currently holds the shared string "Swiss travels" -> "swiss_travels"
$ws.Range("A2").Value = "swiss_travels"

# Author's last selection before saving moved from A2 to A3
$ws.Range("A3").Select()
